$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "tip" text shown for both B14 and B21 used to share the same string
# ("caishuijing"). Split it: B14 keeps the (renamed) original tip, B21 gets
# a brand new tip string.
$ws.Range("B14").Value = "zyshi5"
$ws.Range("B21").Value = "zyyu5"

# Move the active selection to B21, matching where the edit was made.
$ws.Range("B21").Select() | Out-Null
